$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the BOM rows to reflect the Core64 LM v0.3 release values
$ws.Range("A3").Value = "WS2813C"
$ws.Range("D3").Value = "C194323"
$ws.Range("D2").Value = "C14663"

# Update the active selection to A3 (matches the saved selection in the file)
$ws.Range("A3").Select()
